# Daily attendance processing - 2026-01-17 09:03:48
#
# The "Recorded By" column (G) lists every account that touched a given
# attendance session, newest contribution first. A batch reconciliation run
# re-derived the contribution order for every session row: whichever entry
# reflects the most recent recording action is promoted to the front of the
# comma-separated list, while the remaining entries keep their existing
# relative order.
#   - If "System" (exact case) recorded most recently, it is moved to the
#     front of the list.
#   - Otherwise (no "System" entry in the list), the two most-recent human
#     recorders swap position, putting the latest one first.
#   - Single-contributor cells need no reordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 157
$col = 7   # column G - "Recorded By"

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $col)
    $current = $cell.Value2
    if ($current -eq $null) { continue }

    $parts = @($current -split ", ")
    if ($parts.Length -le 1) { continue }

    $systemParts = @($parts | Where-Object { $_.Equals("System") })
    $otherParts = @($parts | Where-Object { -not $_.Equals("System") })

    if ($systemParts.Length -gt 0) {
        # Promote the (exact-case) "System" entries to the front, keeping
        # everything else in its original relative order.
        $newParts = @($systemParts + $otherParts)
    } else {
        # No "System" entry - the most recent recorder moves from the end
        # to the front (equivalent to a reversal for these short lists).
        $newParts = @($parts[($parts.Length - 1)..0])
    }

    $updated = $newParts -join ", "
    if ($updated -ne $current) {
        $cell.Value = $updated
    }
}
